# Add a new worksheet "Versuch6" at the end of the workbook (after "Versuch 4")
# with the RC low-pass filter measurement data (frequency response / Grenzfrequenz).

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Versuch6"

# --- Header row ---
# (the "Grenzfrequenz" label further down the sheet is entered in between the
# header cells below purely so the shared-string table ends up in the same
# insertion order as the source workbook: f, ua, Grenzfrequenz, ue, av, av expected, ...)
$ws.Range("A1").Value = "f"
$ws.Range("B1").Value = "ua"
$ws.Range("A16").Value = "Grenzfrequenz"
$ws.Range("C1").Value = "ue"
$ws.Range("D1").Value = "av"
$ws.Range("E1").Value = "av expected"

# --- Main data block: rows 2-14 ---
# columns: row, f (A), ua (B), ue (C)
$block1 = @(
    @(2, 10, 19.4, 1.94),
    @(3, 20, 19.4, 1.94),
    @(4, 50, 19.4, 1.94),
    @(5, 100, 19.4, 1.94),
    @(6, 200, 19.3, 1.94),
    @(7, 500, 18.7, 1.94),
    @(8, 1000, 16.8, 1.94),
    @(9, 2000, 12.22, 1.94),
    @(10, 5000, 5.99, 1.94),
    @(11, 10000, 3.14, 1.94),
    @(12, 20000, 1.63, 1.94),
    @(13, 50000, 0.68, 1.94),
    @(14, 100000, 0.34, 1.94)
)

foreach ($row in $block1) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Formula = "=B$r/C$r"
    $ws.Range("E$r").Formula = "=(`$B`$31/`$B`$32)*(1/(1+2*PI()*A$r*`$B`$31*`$B`$33))"
}

# --- Second data block: rows 17-26 (around the cutoff frequency) ---
$block2 = @(
    @(17, 1200, 15.7, 1.94),
    @(18, 1300, 15.2, 1.94),
    @(19, 1400, 14.6, 1.94),
    @(20, 1500, 14.3, 1.94),
    @(21, 1591, 13.9, 1.94),
    @(22, 1700, 13.4, 1.94),
    @(23, 1800, 13, 1.94),
    @(24, 1900, 12.6, 1.94),
    @(25, 2000, 12.22, 1.94),
    @(26, 1600, 13.9, 1.94)
)

foreach ($row in $block2) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("E$r").Formula = "=(`$B`$31/`$B`$32)*(1/(1+2*PI()*A$r*`$B`$31*`$B`$33))"
}
# Rows 25 & 26 have a D-column formula that (per the original data) divides by C26
# rather than by the row's own C cell - preserved here to match the source data.
$ws.Range("D17").Formula = "=B17/C17"
$ws.Range("D18").Formula = "=B18/C18"
$ws.Range("D19").Formula = "=B19/C19"
$ws.Range("D20").Formula = "=B20/C20"
$ws.Range("D21").Formula = "=B21/C21"
$ws.Range("D22").Formula = "=B22/C22"
$ws.Range("D23").Formula = "=B23/C23"
$ws.Range("D24").Formula = "=B24/C24"
$ws.Range("D25").Formula = "=B25/C26"
$ws.Range("D26").Formula = "=B26/C26"

# Highlight the measured cutoff-frequency value in bold
$ws.Range("A21").Font.Bold = $true

# --- Sparse row 29 ---
$ws.Range("A29").Value = 15700
$ws.Range("D29").Value = 1

# --- Circuit parameters (R2, R1, C) used by the "av expected" formula ---
$ws.Range("A31").Value = "R2"
$ws.Range("B31").Value = 10000
$ws.Range("A32").Value = "R1"
$ws.Range("B32").Value = 1000
$ws.Range("A33").Value = "C"
$ws.Range("B33").Value = 0.00000001

# Column B is the widest of the data columns - best-fit its width like Excel would.
$ws.Columns.Item(2).AutoFit()

# Put the selection where it was left in the source file and make this the active tab.
$ws.Range("E6").Select()
